$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 37 (shifts existing rows 37-99 down to 38-100,
# and the sheet dimension grows from R99 to R100 automatically).
$ws.Range("A37").EntireRow.Insert()

# Populate the newly inserted row 37 with the new record.
$ws.Range("A37").Value = 6
$ws.Range("B37").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C37").Value = 'Metropolitana'
$ws.Range("D37").Value = 45128
$ws.Range("E37").Value = 13
$ws.Range("F37").Value = 100112035
$ws.Range("G37").Value = 'Bruselas (repollito)'
$ws.Range("H37").Value = 'Sin especificar'
$ws.Range("I37").Value = 'Primera'
$ws.Range("J37").Value = 480
$ws.Range("K37").Value = 17000
$ws.Range("L37").Value = 18000
$ws.Range("M37").Value = 17521
$ws.Range("N37").Value = '$/malla 15 kilos'
$ws.Range("O37").Value = 'Provincia de Quillota'
$ws.Range("P37").Value = 1168
$ws.Range("Q37").Value = 15
$ws.Range("R37").Value = 'Hortaliza'
